$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# A new colo row (ABQ / Albuquerque) was inserted right before the
# existing "ADL" row, so every row from the old 286 ("ADL") through
# the old 297 ("PPT") shifts down by one (to 287..298).
$ws.Rows.Item(286).Insert()

# Copy the A-column cell style (bold/bordered/centered "colo" style)
# from the neighbouring row down onto the new row's A cell so it
# matches every other code cell in the column.
$ws.Cells.Item(287, 1).Copy()
$ws.Cells.Item(286, 1).PasteSpecial(-4122)

# Populate the new row with the Albuquerque colo data.
$ws.Cells.Item(286, 1).Value = "ABQ"
$ws.Cells.Item(286, 2).Value = "Albuquerque"
$ws.Cells.Item(286, 3).Value = 35.0844
$ws.Cells.Item(286, 4).Value = -106.6504
$ws.Cells.Item(286, 5).Value = "US"
$ws.Cells.Item(286, 6).Value = "North America"
$ws.Cells.Item(286, 7).Value = "Albuquerque"
